$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.6
$ws.Range("F4").Value = 0.077
$ws.Range("G4").Value = 0.278
$ws.Range("N4").Value = 0.476
$ws.Range("O4").Value = 0.059
$ws.Range("P4").Value = 0.243
$ws.Range("Q4").Value = 0.19
$ws.Range("R4").Value = 0.122
$ws.Range("S4").Value = 0.35
$ws.Range("W4").Value = 0.381
$ws.Range("X4").Value = 0.06900000000000001
$ws.Range("Y4").Value = 0.263
$ws.Range("AI4").Value = 0.348
$ws.Range("AJ4").Value = 0.102
$ws.Range("AK4").Value = 0.32
$ws.Range("AU4").Value = 0.204
$ws.Range("AV4").Value = 0.036
$ws.Range("AW4").Value = 0.19
$ws.Range("BA4").Value = 1.655
$ws.Range("BB4").Value = 0.079
$ws.Range("BC4").Value = 0.281
$ws.Range("BG4").Value = 0.643
$ws.Range("BH4").Value = 0.194
$ws.Range("BI4").Value = 0.44
$ws.Range("BM4").Value = 0.548
$ws.Range("BN4").Value = 0.045
$ws.Range("BO4").Value = 0.213
$ws.Range("BP4").Value = 0.552
$ws.Range("BQ4").Value = 0.628
$ws.Range("E5").Value = 0.664
$ws.Range("F5").Value = 0.063
$ws.Range("G5").Value = 0.25
$ws.Range("N5").Value = 0.764
$ws.Range("O5").Value = 0.056
$ws.Range("P5").Value = 0.236
$ws.Range("Q5").Value = 0.076
$ws.Range("R5").Value = 0.016
$ws.Range("S5").Value = 0.126
$ws.Range("W5").Value = 0.338
$ws.Range("X5").Value = 0.102
$ws.Range("Y5").Value = 0.319
$ws.Range("AI5").Value = 0.367
$ws.Range("AJ5").Value = 0.099
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.371
$ws.Range("AV5").Value = 0.136
$ws.Range("AW5").Value = 0.368
$ws.Range("BA5").Value = 1.184
$ws.Range("BB5").Value = 0.027
$ws.Range("BC5").Value = 0.165
$ws.Range("BG5").Value = 0.419
$ws.Range("BH5").Value = 0.074
$ws.Range("BI5").Value = 0.272
$ws.Range("BM5").Value = 0.398
$ws.Range("BN5").Value = 0.023
$ws.Range("BO5").Value = 0.152
$ws.Range("BQ5").Value = 0.395
$ws.Range("E6").Value = 0.63
$ws.Range("N6").Value = 0.587
$ws.Range("Q6").Value = 0.109
$ws.Range("W6").Value = 0.358
$ws.Range("AI6").Value = 0.357
$ws.Range("AU6").Value = 0.263
$ws.Range("BA6").Value = 1.378
$ws.Range("BG6").Value = 0.507
$ws.Range("BM6").Value = 0.461
$ws.Range("BP6").Value = 0.459
$ws.Range("BQ6").Value = 0.483
$ws.Range("E7").Value = 0.65
$ws.Range("N7").Value = 0.6820000000000001
$ws.Range("Q7").Value = 0.08599999999999999
$ws.Range("W7").Value = 0.346
$ws.Range("AI7").Value = 0.363
$ws.Range("AU7").Value = 0.319
$ws.Range("BA7").Value = 1.254
$ws.Range("BG7").Value = 0.45
$ws.Range("BM7").Value = 0.421
$ws.Range("BP7").Value = 0.418
$ws.Range("BQ7").Value = 0.426
$ws.Range("E8").Value = 0.633
$ws.Range("F8").Value = 0.08799999999999999
$ws.Range("G8").Value = 0.297
$ws.Range("N8").Value = 0.93
$ws.Range("O8").Value = 0.006
$ws.Range("P8").Value = 0.08
$ws.Range("Q8").Value = 0.112
$ws.Range("R8").Value = 0.055
$ws.Range("S8").Value = 0.234
$ws.Range("W8").Value = 0.475
$ws.Range("X8").Value = 0.108
$ws.Range("Y8").Value = 0.329
$ws.Range("AI8").Value = 0.407
$ws.Range("AJ8").Value = 0.143
$ws.Range("AK8").Value = 0.378
$ws.Range("AU8").Value = 0.31
$ws.Range("AV8").Value = 0.116
$ws.Range("AW8").Value = 0.341
$ws.Range("BA8").Value = 1.657
$ws.Range("BB8").Value = 0.099
$ws.Range("BC8").Value = 0.315
$ws.Range("BG8").Value = 0.538
$ws.Range("BH8").Value = 0.159
$ws.Range("BI8").Value = 0.399
$ws.Range("BM8").Value = 0.638
$ws.Range("BN8").Value = 0.054
$ws.Range("BO8").Value = 0.231
$ws.Range("BP8").Value = 0.552
$ws.Range("BQ8").Value = 0.541
$ws.Range("E9").Value = 0.429
$ws.Range("F9").Value = 0.245
$ws.Range("G9").Value = 0.495
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("AI9").Value = 0.429
$ws.Range("AJ9").Value = 0.245
$ws.Range("AK9").Value = 0.495
$ws.Range("BA9").Value = 1.714
$ws.Range("BM9").Value = 0.714
$ws.Range("BN9").Value = 0.204
$ws.Range("BO9").Value = 0.452
$ws.Range("BP9").Value = 0.571
$ws.Range("BQ9").Value = 0.515
$ws.Range("E10").Value = 0.571
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("W10").Value = 0.571
$ws.Range("AI10").Value = 0.429
$ws.Range("AJ10").Value = 0.245
$ws.Range("AK10").Value = 0.495
$ws.Range("BA10").Value = 1.999
$ws.Range("BM10").Value = 0.857
$ws.Range("BN10").Value = 0.122
$ws.Range("BO10").Value = 0.35
$ws.Range("BP10").Value = 0.666
$ws.Range("BQ10").Value = 0.667
$ws.Range("E11").Value = 0.571
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("W11").Value = 0.571
$ws.Range("AI11").Value = 0.429
$ws.Range("AJ11").Value = 0.245
$ws.Range("AK11").Value = 0.495
$ws.Range("AU11").Value = 0.286
$ws.Range("AV11").Value = 0.204
$ws.Range("AW11").Value = 0.452
$ws.Range("BA11").Value = 1.999
$ws.Range("BM11").Value = 0.857
$ws.Range("BN11").Value = 0.122
$ws.Range("BO11").Value = 0.35
$ws.Range("BP11").Value = 0.666
$ws.Range("BQ11").Value = 0.667
$ws.Range("E12").Value = 1.5
$ws.Range("F12").Value = 0.75
$ws.Range("G12").Value = 0.866
$ws.Range("W12").Value = 1.25
$ws.Range("X12").Value = 0.188
$ws.Range("Y12").Value = 0.433
$ws.Range("AV12").Value = 6
$ws.Range("AW12").Value = 2.449
$ws.Range("BA12").Value = 3.417
$ws.Range("BB12").Value = 0.188
$ws.Range("BC12").Value = 0.433
$ws.Range("BM12").Value = 1.167
$ws.Range("BN12").Value = 0.139
$ws.Range("BO12").Value = 0.373
$ws.Range("BP12").Value = 1.139
$ws.Range("BQ12").Value = 1.279
$ws.Range("BP13").Value = 0.754
$ws.Range("BQ13").Value = 0.6840000000000001
